$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Drop the old "total_lessons_amount" column (T) entirely - it's being
# replaced by per-class teacher/lesson-count pairs in columns D:G.
$ws1.Columns.Item(20).Delete()

# Each class (rows 2-6) now lists three teacher/lesson-count pairs
# (columns B:C, D:E, F:G). Column pair D:E was already partly filled for
# row 2; fill in the rest, then add the new F:G pair for every class.
$data = @(
  @("Brown", 2, "Black", 1),
  @("Brown", 2, "Black", 1),
  @("Brown", 2, "Black", 1),
  @("Brown", 2, "Black", 1),
  @("Brown", 2, "Black", 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $ws1.Range("D$row").Value = $data[$i][0]
  $ws1.Range("E$row").Value = $data[$i][1]
  $ws1.Range("F$row").Value = $data[$i][2]
  $ws1.Range("G$row").Value = $data[$i][3]
}

# Zoom way in on the (now much narrower) table and move the selection
# off the data, matching where the author left the cursor.
$ws1.Application.ActiveWindow.Zoom = 205
$ws1.Range("H3").Select()
